$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 3
$ws.Range("G3").Value = 1.87
$ws.Range("H3").Value = 3.75
$ws.Range("I3").Value = 3.45
$ws.Range("J3").Value = 2.4
$ws.Range("K3").Value = 2.27
$ws.Range("L3").Value = 3.8
$ws.Range("X3").Value = 2.07
$ws.Range("Y3").Value = 8.75
$ws.Range("Z3").Value = 9.75
$ws.Range("AA3").Value = 8.5
$ws.Range("AB3").Value = 16
$ws.Range("AC3").Value = 14
$ws.Range("AH3").Value = 50
$ws.Range("AI3").Value = 12.5
$ws.Range("AJ3").Value = 20
$ws.Range("AK3").Value = 12
$ws.Range("AL3").Value = 45
$ws.Range("AM3").Value = 28

# Row 4
$ws.Range("O4").Value = 1.2
$ws.Range("P4").Value = 4.33
$ws.Range("Q4").Value = 1.7
$ws.Range("W4").Value = 1.73

# Row 5
$ws.Range("H5").Value = 5.5
$ws.Range("J5").Value = 1.73
$ws.Range("K5").Value = 2.63
$ws.Range("M5").Value = 1.02
$ws.Range("N5").Value = 12
$ws.Range("O5").Value = 1.13
$ws.Range("P5").Value = 5.5
$ws.Range("Q5").Value = 1.44
$ws.Range("R5").Value = 2.63
$ws.Range("U5").Value = 1.22
$ws.Range("V5").Value = 4
$ws.Range("W5").Value = 1.73
$ws.Range("X5").Value = 2
$ws.Range("Y5").Value = 10
$ws.Range("Z5").Value = 8.5
$ws.Range("AD5").Value = 21
$ws.Range("AE5").Value = 21
$ws.Range("AH5").Value = 41
$ws.Range("AK5").Value = 21
$ws.Range("AM5").Value = 41

# Row 6
$ws.Range("H6").Value = 5.75
$ws.Range("J6").Value = 1.67
$ws.Range("K6").Value = 2.63
$ws.Range("M6").Value = 1.02
$ws.Range("N6").Value = 11
$ws.Range("O6").Value = 1.14
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 2.5
$ws.Range("S6").Value = 2.1
$ws.Range("T6").Value = 1.67
$ws.Range("U6").Value = 1.22
$ws.Range("V6").Value = 3.75
$ws.Range("W6").Value = 1.83
$ws.Range("X6").Value = 1.83
$ws.Range("Y6").Value = 9.5
$ws.Range("Z6").Value = 7.5
$ws.Range("AA6").Value = 9.5
$ws.Range("AD6").Value = 23
$ws.Range("AE6").Value = 19
$ws.Range("AO6").Value = 600
$ws.Range("AP6").Value = 1.8
$ws.Range("AQ6").Value = 2

# Row 7
$ws.Range("G7").Value = 2.3
$ws.Range("H7").Value = 3.35
$ws.Range("J7").Value = 2.87
$ws.Range("K7").Value = 2.12
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 6.9
$ws.Range("O7").Value = 1.35
$ws.Range("P7").Value = 2.95
$ws.Range("Q7").Value = 2.02
$ws.Range("R7").Value = 1.7
$ws.Range("S7").Value = 3.45
$ws.Range("T7").Value = 1.27
$ws.Range("U7").Value = 1.4
$ws.Range("V7").Value = 2.75
$ws.Range("W7").Value = 1.83
$ws.Range("X7").Value = 1.87
$ws.Range("Y7").Value = 7.3
$ws.Range("AC7").Value = 19.5
$ws.Range("AD7").Value = 32
$ws.Range("AE7").Value = 6.9
$ws.Range("AF7").Value = 6.4
$ws.Range("AG7").Value = 15
$ws.Range("AH7").Value = 75
$ws.Range("AI7").Value = 8.5
$ws.Range("AJ7").Value = 14
$ws.Range("AN7").Value = 37
$ws.Range("AO7").Value = 700
